$wb = $excel.ActiveWorkbook

# Sheets referenced by name (order in workbook: NewProject, kk, NewTask, ProjectTeam, TaskStatus)
$wsNewProject = $wb.Worksheets.Item("NewProject")
$wsNewTask    = $wb.Worksheets.Item("NewTask")
$wsTaskStatus = $wb.Worksheets.Item("TaskStatus")

# 1) NewTask!F5 : "New Task Create Successfully" -> "ForceFully Failed"
#    (done first so the new shared string for it lands before the one created below,
#     matching the recorded shared string ordering)
$wsNewTask.Range("F5").Value = "ForceFully Failed"

# 2) NewProject!A2 and B2 : "POL521" -> "PROJECT001"
$wsNewProject.Range("A2").Value = "PROJECT001"
$wsNewProject.Range("B2").Value = "PROJECT001"

# 3) Update the selected / active cell on each affected sheet.
#    NewProject is selected last so it remains the active (tabSelected) sheet,
#    matching the original workbook where NewProject's sheetView has tabSelected="1".
$wsNewTask.Range("F6").Select()
$wsTaskStatus.Range("B30").Select()
$wsNewProject.Range("E5").Select()
